$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H8").Value = 319.2857
$ws.Range("I8").Value = 206.875
$ws.Range("J8").Value = 335.63635
$ws.Range("K8").Value = 620.625
$ws.Range("L8").Value = 1006.90905
$ws.Range("M8").Value = -481.625
$ws.Range("N8").Value = -1284.90905
$ws.Range("H15").Value = 12324.104
$ws.Range("I15").Value = 12324.104
$ws.Range("K15").Value = 36972.312
$ws.Range("M15").Value = -36803.312
$ws.Range("H18").Value = 15000
$ws.Range("I18").Value = 5000
$ws.Range("J18").Value = 25000
$ws.Range("K18").Value = 5000
$ws.Range("L18").Value = 25000
$ws.Range("M18").Value = -4716
$ws.Range("N18").Value = -25568
$ws.Range("H76").Value = 4357
$ws.Range("I76").Value = 3481.75
$ws.Range("J76").Value = 4857.143
$ws.Range("K76").Value = 3481.75
$ws.Range("L76").Value = 4857.143
$ws.Range("M76").Value = -3166.75
$ws.Range("N76").Value = -5487.143
$ws.Range("H79").Value = 4357
$ws.Range("I79").Value = 3481.75
$ws.Range("J79").Value = 4857.143
$ws.Range("K79").Value = 3481.75
$ws.Range("L79").Value = 4857.143
$ws.Range("M79").Value = -2389.75
$ws.Range("N79").Value = -7041.143
$ws.Range("H113").Value = 6912.6665
$ws.Range("I113").Value = 7199.5557
$ws.Range("K113").Value = 7199.5557
$ws.Range("M113").Value = -3945.5557
$ws.Range("H137").Value = 29241.809
$ws.Range("I137").Value = 15703.8
$ws.Range("K137").Value = 47111.39999999999
$ws.Range("M137").Value = -44561.39999999999
$ws.Range("H138").Value = 21022.727
$ws.Range("I138").Value = 2240.05
$ws.Range("J138").Value = 31755.686
$ws.Range("K138").Value = 6720.150000000001
$ws.Range("L138").Value = 95267.058
$ws.Range("M138").Value = -1580.150000000001
$ws.Range("N138").Value = -105547.058
$ws.Range("H141").Value = 2884.6
$ws.Range("J141").Value = 2800
$ws.Range("L141").Value = 8400
$ws.Range("N141").Value = -18760

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2237.348
$ws.Range("I2").Value = 2154.7144
$ws.Range("J2").Value = 3105
$ws.Range("K2").Value = 2154.7144
$ws.Range("L2").Value = 3105
$ws.Range("M2").Value = -2041.7144
$ws.Range("N2").Value = -3331
$ws.Range("N18").ClearContents()
$ws.Range("H18").Value = 0
$ws.Range("J18").Value = 0
$ws.Range("L18").Value = 0
$ws.Range("H32").Value = 17732.125
$ws.Range("I32").Value = 18341.918
$ws.Range("K32").Value = 18341.918
$ws.Range("M32").Value = -18054.918
$ws.Range("H61").Value = 7510.1577
$ws.Range("I61").Value = 2089.2
$ws.Range("K61").Value = 2089.2
$ws.Range("M61").Value = -1877.2
$ws.Range("H74").Value = 340626.72
$ws.Range("I74").Value = 750800.5
$ws.Range("J74").Value = 12487.7
$ws.Range("K74").Value = 750800.5
$ws.Range("L74").Value = 12487.7
$ws.Range("M74").Value = -749926.5
$ws.Range("N74").Value = -14235.7
$ws.Range("H77").Value = 340626.72
$ws.Range("I77").Value = 750800.5
$ws.Range("J77").Value = 12487.7
$ws.Range("K77").Value = 3754002.5
$ws.Range("L77").Value = 62438.5
$ws.Range("M77").Value = -3749634.5
$ws.Range("N77").Value = -71174.5
$ws.Range("H110").Value = 30623.379
$ws.Range("I110").Value = 34681.52
$ws.Range("K110").Value = 34681.52
$ws.Range("M110").Value = -32636.52
$ws.Range("H116").Value = 2237.348
$ws.Range("I116").Value = 2154.7144
$ws.Range("J116").Value = 3105
$ws.Range("K116").Value = 2154.7144
$ws.Range("L116").Value = 3105
$ws.Range("M116").Value = 139.2856000000002
$ws.Range("N116").Value = -7693
$ws.Range("H132").Value = 1561.075
$ws.Range("I132").Value = 1317.2285
$ws.Range("K132").Value = 3951.6855
$ws.Range("M132").Value = -1421.6855
$ws.Range("H136").Value = 7510.1577
$ws.Range("I136").Value = 2089.2
$ws.Range("K136").Value = 6267.599999999999
$ws.Range("M136").Value = -3717.599999999999
$ws.Range("H139").Value = 114662.664
$ws.Range("J139").Value = 114662.664
$ws.Range("L139").Value = 114662.664
$ws.Range("N139").Value = -124942.664

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2237.348
$ws.Range("I3").Value = 2154.7144
$ws.Range("J3").Value = 3105
$ws.Range("K3").Value = 2154.7144
$ws.Range("L3").Value = 3105
$ws.Range("M3").Value = -2040.7144
$ws.Range("N3").Value = -3333
$ws.Range("H99").Value = 937.9286
$ws.Range("I99").Value = 933.2308
$ws.Range("K99").Value = 933.2308
$ws.Range("M99").Value = 564.7692
$ws.Range("H107").Value = 3009.75
$ws.Range("I107").Value = 3206.389
$ws.Range("K107").Value = 3206.389
$ws.Range("M107").Value = -1286.389
$ws.Range("H132").Value = 89759.664
$ws.Range("J132").Value = 89755
$ws.Range("L132").Value = 89755
$ws.Range("N132").Value = -99875

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1999.6428
$ws.Range("I16").Value = 1768.8462
$ws.Range("J16").Value = 5000
$ws.Range("K16").Value = 1768.8462
$ws.Range("L16").Value = 5000
$ws.Range("M16").Value = -1481.8462
$ws.Range("N16").Value = -5574
$ws.Range("H31").Value = 4000876
$ws.Range("J31").Value = 1250
$ws.Range("L31").Value = 1250
$ws.Range("N31").Value = -1840
$ws.Range("H34").Value = 4000876
$ws.Range("J34").Value = 1250
$ws.Range("L34").Value = 1250
$ws.Range("N34").Value = -1654
$ws.Range("H70").Value = 33491.668
$ws.Range("J70").Value = 33491.668
$ws.Range("L70").Value = 33491.668
$ws.Range("N70").Value = -34121.668
$ws.Range("H73").Value = 33491.668
$ws.Range("J73").Value = 33491.668
$ws.Range("L73").Value = 33491.668
$ws.Range("N73").Value = -35675.668
$ws.Range("H80").Value = 34900
$ws.Range("J80").Value = 34900
$ws.Range("L80").Value = 34900
$ws.Range("N80").Value = -37146
$ws.Range("H83").Value = 34900
$ws.Range("J83").Value = 34900
$ws.Range("L83").Value = 104700
$ws.Range("N83").Value = -115932
$ws.Range("H94").Value = 1494.2727
$ws.Range("I94").Value = 1141.75
$ws.Range("K94").Value = 1141.75
$ws.Range("M94").Value = -690.75
$ws.Range("H99").Value = 6339.385
$ws.Range("I99").Value = 5868.2856
$ws.Range("K99").Value = 5868.2856
$ws.Range("M99").Value = -4370.2856
$ws.Range("H113").Value = 1999.6428
$ws.Range("I113").Value = 1768.8462
$ws.Range("J113").Value = 5000
$ws.Range("K113").Value = 1768.8462
$ws.Range("L113").Value = 5000
$ws.Range("M113").Value = 401.1538
$ws.Range("N113").Value = -9340
$ws.Range("N122").ClearContents()
$ws.Range("H122").Value = 1542.3334
$ws.Range("I122").Value = 1542.3334
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 4627.0002
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -2177.0002
$ws.Range("H126").Value = 6339.385
$ws.Range("I126").Value = 5868.2856
$ws.Range("K126").Value = 17604.8568
$ws.Range("M126").Value = -15134.8568
$ws.Range("H132").Value = 51753.7
$ws.Range("I132").Value = 67560.60000000001
$ws.Range("K132").Value = 202681.8
$ws.Range("M132").Value = -200151.8
$ws.Range("H134").Value = 3423
$ws.Range("I134").Value = 2661.125
$ws.Range("J134").Value = 4438.8335
$ws.Range("K134").Value = 7983.375
$ws.Range("L134").Value = 13316.5005
$ws.Range("M134").Value = -5448.375
$ws.Range("N134").Value = -18386.5005

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 1024.4286
$ws.Range("J122").Value = 1126.5454
$ws.Range("L122").Value = 10138.9086
$ws.Range("N122").Value = -15038.9086
$ws.Range("H137").Value = 2795.375
$ws.Range("I137").Value = 2402.7693
$ws.Range("J137").Value = 4496.6665
$ws.Range("K137").Value = 7208.3079
$ws.Range("L137").Value = 13489.9995
$ws.Range("M137").Value = -2108.3079
$ws.Range("N137").Value = -23689.9995

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H41").Value = 8779.4
$ws.Range("I41").Value = 1998.5
$ws.Range("J41").Value = 13300
$ws.Range("K41").Value = 1998.5
$ws.Range("L41").Value = 13300
$ws.Range("M41").Value = -1643.5
$ws.Range("N41").Value = -14010
$ws.Range("N101").ClearContents()
$ws.Range("H101").Value = 0
$ws.Range("J101").Value = 0
$ws.Range("L101").Value = 0
$ws.Range("H107").Value = 733.3
$ws.Range("I107").Value = 735.1429000000001
$ws.Range("K107").Value = 735.1429000000001
$ws.Range("M107").Value = 1184.8571
$ws.Range("H132").Value = 2327.8276
$ws.Range("I132").Value = 1772.875
$ws.Range("J132").Value = 4991.6
$ws.Range("K132").Value = 5318.625
$ws.Range("L132").Value = 14974.8
$ws.Range("M132").Value = -2788.625
$ws.Range("N132").Value = -20034.8
$ws.Range("H140").Value = 30000
$ws.Range("J140").Value = 80000
$ws.Range("L140").Value = 80000
$ws.Range("N140").Value = -90360

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 5998.4
$ws.Range("I136").Value = 5831.3335
$ws.Range("K136").Value = 17494.0005
$ws.Range("M136").Value = -14944.0005

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H14").Value = 5750
$ws.Range("I14").Value = 5000
$ws.Range("K14").Value = 5000
$ws.Range("M14").Value = -4832
$ws.Range("H45").Value = 17782.584
$ws.Range("J45").Value = 17782.584
$ws.Range("L45").Value = 17782.584
$ws.Range("N45").Value = -18764.584
$ws.Range("H132").Value = 46870.062
$ws.Range("I132").Value = 56994.31
$ws.Range("J132").Value = 2998.3333
$ws.Range("K132").Value = 170982.93
$ws.Range("L132").Value = 8994.999899999999
$ws.Range("M132").Value = -168452.93
$ws.Range("N132").Value = -14054.9999
$ws.Range("H136").Value = 25513.74
$ws.Range("I136").Value = 28265.459
$ws.Range("K136").Value = 84796.37699999999
$ws.Range("M136").Value = -82246.37699999999
